$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.030.64'
$ws.Range("E2").Value = '  -3.37%  '

$ws.Range("D3").Value = '2.897.05'
$ws.Range("E3").Value = '  -3.28%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = "'369.49"
$ws.Range("E5").Value = '  +2.81%  '

$ws.Range("D6").Value = "'102.63"
$ws.Range("E6").Value = '  -7.26%  '

$ws.Range("E7").Value = '  -6.06%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = '  -7.75%  '

$ws.Range("D10").Value = "'36.82"
$ws.Range("E10").Value = '  -6.99%  '

$ws.Range("E11").Value = '  +0.81%  '

$ws.Range("D12").Value = "'0.0835"
$ws.Range("E12").Value = '  -5.57%  '

$ws.Range("D13").Value = "'18.41"
$ws.Range("E13").Value = '  -5.93%  '

$ws.Range("D14").Value = '3.354.64'
$ws.Range("E14").Value = '  -2.98%  '

$ws.Range("E15").Value = '  -6.71%  '

$ws.Range("D16").Value = '2.896.04'
$ws.Range("E16").Value = '  -3.66%  '

$ws.Range("E17").Value = '  -6.39%  '

$ws.Range("D18").Value = '50.939.02'
$ws.Range("E18").Value = '  -3.51%  '

$ws.Range("D19").Value = "'3.28"
$ws.Range("E19").Value = '  -6.00%  '

$ws.Range("E20").Value = '  -6.02%  '

$ws.Range("D21").Value = "'12.90"
$ws.Range("E21").Value = '  -8.20%  '

$ws.Range("D22").Value = '0.0₃0943'
$ws.Range("E22").Value = '  -5.11%  '

$ws.Range("D23").Value = "'68.26"
$ws.Range("E23").Value = '  -3.70%  '

$ws.Range("D24").Value = "'259.56"
$ws.Range("E24").Value = '  -5.00%  '

$ws.Range("E25").Value = '  -4.49%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = "'0.169"
$ws.Range("E26").Value = '  -6.92%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = '  -0.13%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'25.77"
$ws.Range("E28").Value = '  -6.25%  '

$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = "'7.12"
$ws.Range("E29").Value = '  -7.24%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = "'0.104"
$ws.Range("E30").Value = '  -5.29%  '

$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = "'9.89"
$ws.Range("E31").Value = '  -6.65%  '

$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D32").Value = "'6.03"
$ws.Range("E32").Value = '  -1.77%  '

$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").Value = "'2.13"
$ws.Range("E33").Value = '  -2.85%  '

$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = "'34.96"
$ws.Range("E34").Value = '  -8.84%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = "'51.04"
$ws.Range("E35").Value = '  -2.69%  '

$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = '  +0.41%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.0418"
$ws.Range("E37").Value = '  -6.80%  '

$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = "'3.06"
$ws.Range("E38").Value = '  -6.60%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = "'2.66"
$ws.Range("E39").Value = '  -5.17%  '

$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = "'16.92"
$ws.Range("E40").Value = '  -7.41%  '

$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = "'1.84"
$ws.Range("E41").Value = '  -10.27%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = "'0.113"
$ws.Range("E42").Value = '  -6.55%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = "'22.06"
$ws.Range("E43").Value = '  -7.91%  '

$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").Value = "'117.29"
$ws.Range("E44").Value = '  -2.11%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'2.10"
$ws.Range("E45").Value = '  -3.33%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.033.83'
$ws.Range("E46").Value = '  -5.77%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = "'2.32"
$ws.Range("E47").Value = '  -6.28%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'3.17"
$ws.Range("E48").Value = '  -9.50%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '3.189.98'
$ws.Range("E49").Value = '  -2.72%  '

$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = "'0.239"
$ws.Range("E50").Value = '  -3.21%  '

$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").Value = "'0.0312"
$ws.Range("E51").Value = '  -12.35%  '
